$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = 4
$ws.Range("C15").Value = "Count Negative Numbers in a Sorted Matrix"
$ws.Range("D15").Value = "LeetCode"

$ws.Range("C16").Value = "Squares of a Sorted Array"
$ws.Range("D16").Value = "Bosscoder Academy"

$ws.Columns.Item(4).ColumnWidth = 15.67

$ws.Range("E16").Select()
